# Final dataset + final decoding scripts
# Clears the "x" marker from the "LSB Classic / Squares" (E) and
# "Invoke-PSImage" (H) columns (plus one stray F21) on the "Images" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Images")

$cellsToClear = @(
    "E3", "E4", "H5", "E6", "E7", "E8", "H8", "E9", "E10", "H10",
    "E11", "H11", "E12", "E13", "E14", "E15", "E16", "H18", "E19", "H19",
    "E20", "H20", "F21", "H21", "E22", "E23", "E24", "H24", "H25", "E26",
    "E27", "E28", "H28", "H29", "H30", "H31", "H32", "H33", "H34", "E35",
    "H35", "E36", "H36", "E37", "H37", "E38", "H38"
)

foreach ($addr in $cellsToClear) {
    $ws.Range($addr).Value = $null
}

# Move the active selection to match the saved view state.
$ws.Range("E18").Select() | Out-Null
